$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2162.4546
$ws.Range("I8").Value = 22.777779
$ws.Range("J8").Value = 3643.7693
$ws.Range("K8").Value = 68.333337
$ws.Range("L8").Value = 10931.3079
$ws.Range("M8").Value = 70.666663
$ws.Range("N8").Value = -11209.3079

$ws.Range("H19").Value = 1864.3334
$ws.Range("I19").Value = 2097.25
$ws.Range("J19").Value = 1398.5
$ws.Range("K19").Value = 2097.25
$ws.Range("L19").Value = 1398.5
$ws.Range("M19").Value = -1922.25
$ws.Range("N19").Value = -1748.5

$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -532

$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -766

$ws.Range("H33").Value = 150.58333
$ws.Range("I33").Value = 150.58333
$ws.Range("K33").Value = 150.58333
$ws.Range("M33").Value = 78.41667000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 632.3
$ws.Range("I2").Value = 526.5
$ws.Range("K2").Value = 526.5
$ws.Range("M2").Value = -413.5

$ws.Range("H102").Value = 2749.5
$ws.Range("I102").Value = 2749.5
$ws.Range("K102").Value = 2749.5
$ws.Range("M102").Value = -1127.5

$ws.Range("H116").Value = 632.3
$ws.Range("I116").Value = 526.5
$ws.Range("K116").Value = 526.5
$ws.Range("M116").Value = 1767.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 632.3
$ws.Range("I3").Value = 526.5
$ws.Range("K3").Value = 526.5
$ws.Range("M3").Value = -412.5

$ws.Range("H105").Value = 2314.1667
$ws.Range("I105").Value = 1721.25
$ws.Range("K105").Value = 1721.25
$ws.Range("M105").Value = 25.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1858
$ws.Range("J58").Value = 2482.3333
$ws.Range("L58").Value = 2482.3333
$ws.Range("N58").Value = -2888.3333

$ws.Range("H86").Value = 11014.857
$ws.Range("I86").Value = 11853
$ws.Range("K86").Value = 11853
$ws.Range("M86").Value = -10730

$ws.Range("H89").Value = 11014.857
$ws.Range("I89").Value = 11853
$ws.Range("K89").Value = 59265
$ws.Range("M89").Value = -53649

$ws.Range("H99").Value = 1676
$ws.Range("I99").Value = 1419
$ws.Range("J99").Value = 2006.4286
$ws.Range("K99").Value = 1419
$ws.Range("L99").Value = 2006.4286
$ws.Range("M99").Value = 79
$ws.Range("N99").Value = -5002.4286

$ws.Range("H126").Value = 1676
$ws.Range("I126").Value = 1419
$ws.Range("J126").Value = 2006.4286
$ws.Range("K126").Value = 4257
$ws.Range("L126").Value = 6019.2858
$ws.Range("M126").Value = -1787
$ws.Range("N126").Value = -10959.2858

$ws.Range("H136").Value = 1858
$ws.Range("J136").Value = 2482.3333
$ws.Range("L136").Value = 7446.999899999999
$ws.Range("N136").Value = -12546.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H129").Value = 2171.4614
$ws.Range("I129").Value = 1587.3334
$ws.Range("J129").Value = 2672.1428
$ws.Range("K129").Value = 4762.0002
$ws.Range("L129").Value = 8016.428400000001
$ws.Range("M129").Value = 237.9997999999996
$ws.Range("N129").Value = -18016.4284

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

$ws.Range("H131").Value = 2440.2104
$ws.Range("I131").Value = 1483
$ws.Range("J131").Value = 2998.5833
$ws.Range("K131").Value = 4449
$ws.Range("L131").Value = 8995.749899999999
$ws.Range("M131").Value = 591
$ws.Range("N131").Value = -19075.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 88050.664
$ws.Range("J10").Value = 3434.6667
$ws.Range("L10").Value = 3434.6667
$ws.Range("N10").Value = -3772.6667

$ws.Range("H70").Value = 2452
$ws.Range("I70").Value = 2452
$ws.Range("K70").Value = 2452
$ws.Range("M70").Value = -2182

$ws.Range("H73").Value = 2452
$ws.Range("I73").Value = 2452
$ws.Range("K73").Value = 2452
$ws.Range("M73").Value = -1516

$ws.Range("H80").Value = 2871.75
$ws.Range("J80").Value = 3199.75
$ws.Range("L80").Value = 3199.75
$ws.Range("N80").Value = -5195.75

$ws.Range("H83").Value = 2871.75
$ws.Range("J83").Value = 3199.75
$ws.Range("L83").Value = 15998.75
$ws.Range("N83").Value = -25982.75

$ws.Range("H101").Value = 37990
$ws.Range("J101").Value = 37990
$ws.Range("L101").Value = 37990
$ws.Range("N101").Value = -44480

$ws.Range("H126").Value = 2399.8
$ws.Range("I126").Value = 1999.75
$ws.Range("K126").Value = 5999.25
$ws.Range("M126").Value = -3529.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 501.5
$ws.Range("I12").Value = 999
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = 999
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = -829
$ws.Range("N12").Value = -344

$ws.Range("H46").Value = 1671.5555
$ws.Range("I46").Value = 1049.6666
$ws.Range("K46").Value = 1049.6666
$ws.Range("M46").Value = -861.6666

$ws.Range("H55").Value = 1173.5714
$ws.Range("I55").Value = 329.25
$ws.Range("K55").Value = 329.25
$ws.Range("M55").Value = -156.25

$ws.Range("H61").Value = 1349.8334
$ws.Range("I61").Value = 1339.9
$ws.Range("J61").Value = 1399.5
$ws.Range("K61").Value = 1339.9
$ws.Range("L61").Value = 1399.5
$ws.Range("M61").Value = -1137.9
$ws.Range("N61").Value = -1803.5

$ws.Range("H82").Value = 1198.8334
$ws.Range("I82").Value = 798.5
$ws.Range("J82").Value = 1399
$ws.Range("K82").Value = 798.5
$ws.Range("L82").Value = 1399
$ws.Range("M82").Value = -437.5
$ws.Range("N82").Value = -2121

$ws.Range("H85").Value = 1198.8334
$ws.Range("I85").Value = 798.5
$ws.Range("J85").Value = 1399
$ws.Range("K85").Value = 798.5
$ws.Range("L85").Value = 1399
$ws.Range("M85").Value = 449.5
$ws.Range("N85").Value = -3895

$ws.Range("H93").Value = 1700
$ws.Range("I93").Value = 1650
$ws.Range("J93").Value = 1750
$ws.Range("K93").Value = 1650
$ws.Range("L93").Value = 1750
$ws.Range("M93").Value = -402
$ws.Range("N93").Value = -4246

$ws.Range("H113").Value = 1349.8334
$ws.Range("I113").Value = 1339.9
$ws.Range("J113").Value = 1399.5
$ws.Range("K113").Value = 1339.9
$ws.Range("L113").Value = 1399.5
$ws.Range("M113").Value = 830.0999999999999
$ws.Range("N113").Value = -5739.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 17750
$ws.Range("J7").Value = 17750
$ws.Range("L7").Value = 17750
$ws.Range("N7").Value = -17976

$ws.Range("H9").Value = 1999.6
$ws.Range("I9").Value = 2502
$ws.Range("K9").Value = 2502
$ws.Range("M9").Value = -2362

$ws.Range("H12").Value = 2466.6667
$ws.Range("J12").Value = 2466.6667
$ws.Range("L12").Value = 2466.6667
$ws.Range("N12").Value = -2750.6667

$ws.Range("H132").Value = 1835.909
$ws.Range("I132").Value = 1835.909
$ws.Range("K132").Value = 5507.727000000001
$ws.Range("M132").Value = -2977.727000000001
